$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PC1 (col B) and PC2 (col C) values for rows 2-17
$values = @(
    @(0.1303827324240858,   0.2489604785350144),
    @(0.1628699522827561,  -0.1164493770259173),
    @(0.1209827266396085,  -0.3855718635361127),
    @(0.205666759012536,    0.2798943549036084),
    @(0.2928904243400641,  -0.402266262391767),
    @(0.3243971146311336,  -0.2808318545561967),
    @(0.4269953608783478,  -0.1636979571852554),
    @(0.3322062125579185,   0.3948871054647279),
    @(-0.08635994341456091, 0.2795858829863248),
    @(-0.3113047488585209, -0.1260370359714243),
    @(-0.125009764327608,  -0.2393817736955996),
    @(0.01995906236735627,  0.2012386853443357),
    @(0.4931326245915859,   0.1745595534604872),
    @(0.09893399808120779, -0.2219910426751996),
    @(0.2041321347028903,   0.04928163274642373),
    @(-0.02831341782287976,-0.0004524403263734084)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
}
